# Automatische test-sync: 2025-06-17 22:27:08
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# Add new row 51: Afmelding nieuwsbrief
$ws.Range("A51").Value = "Afmelding nieuwsbrief"
$ws.Range("B51").Value = "mailmind.test@zohomail.eu"
$ws.Range("C51").Value = "Graag afmelden voor de nieuwsbrief. Dank u."
$ws.Range("D51").Value = "Afmelding"
$ws.Range("F51").Value = "2025-06-17 22:14:31"
$ws.Range("G51").Value = "Nee"

# Add new row 52: Klacht over levering
$ws.Range("A52").Value = "Klacht over levering"
$ws.Range("B52").Value = "mailmind.test@zohomail.eu"
$ws.Range("C52").Value = "Ik ben niet tevreden over mijn bestelling. Ik hoor graag hoe jullie dit oplossen."
$ws.Range("D52").Value = "Klacht"
$ws.Range("F52").Value = "2025-06-17 22:27:02"
$ws.Range("G52").Value = "Nee"

# Update Dashboard counts
$dash.Range("B4").Value = 10
$dash.Range("B5").Value = 5

# Extend conditional formatting ranges to cover the new rows
$catFcs = $ws.Range("D2:D50").FormatConditions
for ($i = 1; $i -le $catFcs.Count; $i++) {
    $catFcs.Item($i).ModifyAppliesToRange($ws.Range("D2:D52"))
}

$answeredFcs = $ws.Range("G2:G50").FormatConditions
for ($i = 1; $i -le $answeredFcs.Count; $i++) {
    $answeredFcs.Item($i).ModifyAppliesToRange($ws.Range("G2:G52"))
}
